# Clear the "A" (absent) marks recorded in column Y (lecture/day 19)
# for all students that currently have one. The dependent "Total
# Absence" formulas in column E recompute automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(8,9,23,24,26,28,32,33,36,40,41,42,47,49,50,51,57,62,64,68,71,79)

foreach ($r in $rows) {
    $cell = $ws.Range("Y$r")
    $cell.ClearContents()

    # Match the resulting blank-cell formatting (border/font, no explicit
    # alignment) used by the rest of the already-empty day columns: copy
    # formats from the neighbouring blank cell in the same row instead of
    # leaving the old "marked absent" style behind.
    $blankTemplate = $ws.Range("Z$r")
    $blankTemplate.Copy()
    $cell.PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
